$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the TOT label first so the new shared string lands at index 8,
# matching the target order (TOT, Interno, Ricerca info, CM, XAMPP).
$ws.Range("H14").Value = "TOT"

# Row 2 - changed to Interno / Ricerca info, new date/minutes
$ws.Range("F2").Value = "Interno"
$ws.Range("G2").Value = "Ricerca info"
$ws.Range("H2").Value = 43423
$ws.Range("I2").Value = 60

# Row 3 - changed to Interno / Ricerca info, new date/minutes
$ws.Range("F3").Value = "Interno"
$ws.Range("G3").Value = "Ricerca info"
$ws.Range("H3").Value = 43424
$ws.Range("I3").Value = 60

# Row 4 - date/minutes change only
$ws.Range("H4").Value = 43426
$ws.Range("I4").Value = 30

# Row 5 - date change only
$ws.Range("H5").Value = 43429

# Row 6 - changed to Interno / CM, new date/minutes
$ws.Range("F6").Value = "Interno"
$ws.Range("G6").Value = "CM"
$ws.Range("H6").Value = 43434
$ws.Range("I6").Value = 60

# Row 7 - date/minutes change only
$ws.Range("H7").Value = 43444
$ws.Range("I7").Value = 60

# Row 8 - brand new row
$ws.Range("E8").Value = "Gianluca"
$ws.Range("F8").Value = "GDPR"
$ws.Range("G8").Value = "Documentazione"
$ws.Range("H8").Value = 43467
$ws.Range("I8").Value = 30

# Row 9 - brand new row
$ws.Range("E9").Value = "Gianluca"
$ws.Range("F9").Value = "GDPR"
$ws.Range("G9").Value = "Documentazione"
$ws.Range("H9").Value = 43481
$ws.Range("I9").Value = 90

# Row 10 - brand new row (XAMPP)
$ws.Range("E10").Value = "Gianluca"
$ws.Range("F10").Value = "GDPR"
$ws.Range("G10").Value = "XAMPP"
$ws.Range("H10").Value = 43497
$ws.Range("I10").Value = 90

# Row 11 - brand new row (XAMPP)
$ws.Range("E11").Value = "Gianluca"
$ws.Range("F11").Value = "GDPR"
$ws.Range("G11").Value = "XAMPP"
$ws.Range("H11").Value = 43511
$ws.Range("I11").Value = 60

# Row 12 - brand new row
$ws.Range("E12").Value = "Gianluca"
$ws.Range("F12").Value = "GDPR"
$ws.Range("G12").Value = "Documentazione"
$ws.Range("H12").Value = 43526
$ws.Range("I12").Value = 20

# Row 13 - brand new row
$ws.Range("E13").Value = "Gianluca"
$ws.Range("F13").Value = "GDPR"
$ws.Range("G13").Value = "Documentazione"
$ws.Range("H13").Value = 43549
$ws.Range("I13").Value = 90

# Row 14 - totals row (label already set above)
$ws.Range("I14").Formula = "=SUM(I2:I13)"

# Rows 15-17 - empty cells carrying the date style (copy format only, no value)
$ws.Range("H2").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New selection matches the post-edit workbook state
$ws.Range("I14").Select()
